# Add other Consensus Economics forecasts
# - inserts forecast_con_*, forecast_inv_*, forecast_ip_* columns between the
#   existing forecast_gdp_* and forecast_inf_* columns
# - the two forecast_inf_* columns move from K:L to the end (P:Q)
# - forecast_gdp_2step / forecast_inf_2step columns are dropped
# - all forecast_* values are refreshed with new (rounded) consensus figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 48

# --- Header row -------------------------------------------------------
$ws.Range("J1").Value = "forecast_con_current"
$ws.Range("K1").Value = "forecast_con_1step"
$ws.Range("L1").Value = "forecast_inv_current"
$ws.Range("M1").Value = "forecast_inv_1step"
$ws.Range("N1").Value = "forecast_ip_current"
$ws.Range("O1").Value = "forecast_ip_1step"
$ws.Range("P1").Value = "forecast_inf_current"
$ws.Range("Q1").Value = "forecast_inf_1step"

# --- Data rows ---------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value  = -3.8    # H forecast_gdp_current
    $ws.Cells.Item($r, 9).Value  = 3.1     # I forecast_gdp_1step
    $ws.Cells.Item($r, 10).Value = -6.5    # J forecast_con_current
    $ws.Cells.Item($r, 11).Value = 3.8     # K forecast_con_1step
    $ws.Cells.Item($r, 12).Value = -7.8    # L forecast_inv_current
    $ws.Cells.Item($r, 13).Value = 4.1     # M forecast_inv_1step
    $ws.Cells.Item($r, 14).Value = -3.5    # N forecast_ip_current
    $ws.Cells.Item($r, 15).Value = 3.2     # O forecast_ip_1step
    $ws.Cells.Item($r, 16).Value = 3.9     # P forecast_inf_current
    $ws.Cells.Item($r, 17).Value = 3.6     # Q forecast_inf_1step
}

# --- Selection / view state --------------------------------------------
$ws.Range("H3:Q48").Select()
